$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Fixed table sort" - a handful of rows had drifted out of the
# WARSCORE (column J) descending sort order. Re-sort those rows by
# swapping the player (name + stat) between the two rows, then apply
# the "last update" refresh to the WARSCORE column.

function Swap-Row($rowA, $rowB) {
    # Swap name (A) and stat columns (B:I) between two rows; column J
    # (WARSCORE) is left alone here and fixed up separately below.
    $colsToSwap = @("A", "B", "C", "D", "E", "F", "G", "H", "I")
    foreach ($col in $colsToSwap) {
        $refA = "$col$rowA"
        $refB = "$col$rowB"
        $valA = $ws.Range($refA).Value2
        $valB = $ws.Range($refB).Value2
        $ws.Range($refA).Value2 = $valB
        $ws.Range($refB).Value2 = $valA
    }
}

# Rows 25/26 (topdosl33ts / qzt) swap places.
Swap-Row 25 26

# Rows 39/40 (Marcelo / Pikaya) swap places.
Swap-Row 39 40

# Rows 45/46 (SELFIE / Sra. Clash) swap places.
Swap-Row 45 46

# "last update" - refreshed WARSCORE (column J) values.
$ws.Range("J9").Value2 = 222
$ws.Range("J23").Value2 = 152
$ws.Range("J25").Value2 = 122
$ws.Range("J26").Value2 = 122
$ws.Range("J27").Value2 = 120
$ws.Range("J29").Value2 = 104
$ws.Range("J39").Value2 = 48
$ws.Range("J40").Value2 = 48
$ws.Range("J44").Value2 = 15
$ws.Range("J45").Value2 = 3
$ws.Range("J46").Value2 = 3
